# Auto-generated Excel COM-interop script implementing the commit
# "Generate Report for Handoff": reorders per-file rows on the
# Overview/zh-cn/de-de sheets (984b07f3... moves from row 2 to row 4
# and is marked "Ready for handoff"), refreshes hyperlink display text,
# and widens the Error Detail column.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# --- Overview sheet: updated cell values ---
$ws1.Range("A2").Value = 'ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md'
$ws1.Range("B2").Value = 'e2e\ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md'
$ws1.Range("G2").Value = '2016-08-16 11:03:35'
$ws1.Range("A3").Value = 'ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md'
$ws1.Range("B3").Value = 'e2e\ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md'
$ws1.Range("A4").Value = '984b07f3-8c45-4712-a89f-06216d9f9206.md'
$ws1.Range("B4").Value = 'e2e\984b07f3-8c45-4712-a89f-06216d9f9206.md'
$ws1.Range("E4").Value = 'Ready for handoff'
$ws1.Range("F4").Value = 'Ready for handoff'
$ws1.Range("G4").Value = '2016-08-16 11:05:48'

# --- zh-cn sheet: updated cell values ---
$ws2.Range("A2").Value = 'ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md'
$ws2.Range("G2").Value = '885e7002-3dba-40ab-a7d3-33d242224785.0271c48b8eceb70fe07976a53047e849de7936c1.zh-cn.xlf'
$ws2.Range("H2").Value = '2016-08-16 11:03:29'
$ws2.Range("I2").Value = '885e7002-3dba-40ab-a7d3-33d242224785.md'
$ws2.Range("J2").Value = '885e7002-3dba-40ab-a7d3-33d242224785.0271c48b8eceb70fe07976a53047e849de7936c1.zh-cn.xlf'
$ws2.Range("K2").Value = '2016-08-16 11:03:57'
$ws2.Range("A3").Value = 'ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md'
$ws2.Range("F3").Value = 'True'
$ws2.Range("A4").Value = '984b07f3-8c45-4712-a89f-06216d9f9206.md'
$ws2.Range("C4").Value = 'Ready for handoff'
$ws2.Range("F4").Value = 'False'
$ws2.Range("G4").Value = '984b07f3-8c45-4712-a89f-06216d9f9206.0a683d6ce457ecb89daf73c135c76f920d7d20cd.zh-cn.xlf'
$ws2.Range("H4").Value = '2016-08-16 11:05:43'
$ws2.Range("I4").Value = '984b07f3-8c45-4712-a89f-06216d9f9206.md'
$ws2.Range("J4").Value = '984b07f3-8c45-4712-a89f-06216d9f9206.0a683d6ce457ecb89daf73c135c76f920d7d20cd.zh-cn.xlf'
$ws2.Range("K4").Value = '2016-08-16 11:05:21'
$ws2.Range("P4").Value = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa1a6738a64140972add146ece41f8f017b924ad/e2e/984b07f3-8c45-4712-a89f-06216d9f9206.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5fd5174a6669efe7a903696233297b56cbea9b70/e2e/984b07f3-8c45-4712-a89f-06216d9f9206.md.'

# --- de-de sheet: updated cell values ---
$ws3.Range("A2").Value = 'ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md'
$ws3.Range("G2").Value = '885e7002-3dba-40ab-a7d3-33d242224785.0271c48b8eceb70fe07976a53047e849de7936c1.de-de.xlf'
$ws3.Range("H2").Value = '2016-08-16 11:03:35'
$ws3.Range("I2").Value = '885e7002-3dba-40ab-a7d3-33d242224785.md'
$ws3.Range("J2").Value = '885e7002-3dba-40ab-a7d3-33d242224785.0271c48b8eceb70fe07976a53047e849de7936c1.de-de.xlf'
$ws3.Range("K2").Value = '2016-08-16 11:04:12'
$ws3.Range("A3").Value = 'ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md'
$ws3.Range("F3").Value = 'True'
$ws3.Range("A4").Value = '984b07f3-8c45-4712-a89f-06216d9f9206.md'
$ws3.Range("C4").Value = 'Ready for handoff'
$ws3.Range("F4").Value = 'False'
$ws3.Range("G4").Value = '984b07f3-8c45-4712-a89f-06216d9f9206.0a683d6ce457ecb89daf73c135c76f920d7d20cd.de-de.xlf'
$ws3.Range("H4").Value = '2016-08-16 11:05:48'
$ws3.Range("I4").Value = '984b07f3-8c45-4712-a89f-06216d9f9206.md'
$ws3.Range("J4").Value = '984b07f3-8c45-4712-a89f-06216d9f9206.0a683d6ce457ecb89daf73c135c76f920d7d20cd.de-de.xlf'
$ws3.Range("K4").Value = '2016-08-16 11:05:29'
$ws3.Range("P4").Value = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa1a6738a64140972add146ece41f8f017b924ad/e2e/984b07f3-8c45-4712-a89f-06216d9f9206.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5fd5174a6669efe7a903696233297b56cbea9b70/e2e/984b07f3-8c45-4712-a89f-06216d9f9206.md.'

# --- Overview sheet: rebuild hyperlinks (same targets, refreshed display text) ---
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa1a6738a64140972add146ece41f8f017b924ad/e2e/984b07f3-8c45-4712-a89f-06216d9f9206.md', [Type]::Missing, [Type]::Missing, 'e2e\ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md')
$ws1.Hyperlinks.Add($ws1.Range("B3"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40ed266e45e0a8e86a73653470b6acf9515cba7e/e2e/ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md', [Type]::Missing, [Type]::Missing, 'e2e\ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md')
$ws1.Hyperlinks.Add($ws1.Range("B4"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa1a6738a64140972add146ece41f8f017b924ad/e2e/ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md', [Type]::Missing, [Type]::Missing, 'e2e\984b07f3-8c45-4712-a89f-06216d9f9206.md')

# --- zh-cn sheet: rebuild hyperlinks (same targets, refreshed display text) ---
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa1a6738a64140972add146ece41f8f017b924ad/e2e/984b07f3-8c45-4712-a89f-06216d9f9206.md', [Type]::Missing, [Type]::Missing, 'ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md')
$ws2.Hyperlinks.Add($ws2.Range("I2"), 'https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/88d679e9e565835f3ddff4bb9651fc036180cedd/e2e/984b07f3-8c45-4712-a89f-06216d9f9206.md', [Type]::Missing, [Type]::Missing, '885e7002-3dba-40ab-a7d3-33d242224785.md')
$ws2.Hyperlinks.Add($ws2.Range("A3"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40ed266e45e0a8e86a73653470b6acf9515cba7e/e2e/ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md', [Type]::Missing, [Type]::Missing, 'ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md')
$ws2.Hyperlinks.Add($ws2.Range("I3"), 'https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/350b13fbd54e8811bb1b70c6ba16f41dce59c601/e2e/885e7002-3dba-40ab-a7d3-33d242224785.md', [Type]::Missing, [Type]::Missing, '885e7002-3dba-40ab-a7d3-33d242224785.md')
$ws2.Hyperlinks.Add($ws2.Range("A4"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa1a6738a64140972add146ece41f8f017b924ad/e2e/ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md', [Type]::Missing, [Type]::Missing, '984b07f3-8c45-4712-a89f-06216d9f9206.md')
$ws2.Hyperlinks.Add($ws2.Range("I4"), 'https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/350b13fbd54e8811bb1b70c6ba16f41dce59c601/e2e/885e7002-3dba-40ab-a7d3-33d242224785.md', [Type]::Missing, [Type]::Missing, '984b07f3-8c45-4712-a89f-06216d9f9206.md')

# --- de-de sheet: rebuild hyperlinks (same targets, refreshed display text) ---
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa1a6738a64140972add146ece41f8f017b924ad/e2e/984b07f3-8c45-4712-a89f-06216d9f9206.md', [Type]::Missing, [Type]::Missing, 'ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md')
$ws3.Hyperlinks.Add($ws3.Range("I2"), 'https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/916179c27494a5c82077e51d3a552f5b2a6a6233/e2e/984b07f3-8c45-4712-a89f-06216d9f9206.md', [Type]::Missing, [Type]::Missing, '885e7002-3dba-40ab-a7d3-33d242224785.md')
$ws3.Hyperlinks.Add($ws3.Range("A3"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40ed266e45e0a8e86a73653470b6acf9515cba7e/e2e/ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md', [Type]::Missing, [Type]::Missing, 'ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md')
$ws3.Hyperlinks.Add($ws3.Range("I3"), 'https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/e1441ac9d6831a21cf715f1dc85131d2ef0c9a1e/e2e/885e7002-3dba-40ab-a7d3-33d242224785.md', [Type]::Missing, [Type]::Missing, '885e7002-3dba-40ab-a7d3-33d242224785.md')
$ws3.Hyperlinks.Add($ws3.Range("A4"), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa1a6738a64140972add146ece41f8f017b924ad/e2e/ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md', [Type]::Missing, [Type]::Missing, '984b07f3-8c45-4712-a89f-06216d9f9206.md')
$ws3.Hyperlinks.Add($ws3.Range("I4"), 'https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/e1441ac9d6831a21cf715f1dc85131d2ef0c9a1e/e2e/885e7002-3dba-40ab-a7d3-33d242224785.md', [Type]::Missing, [Type]::Missing, '984b07f3-8c45-4712-a89f-06216d9f9206.md')

# --- widen the "Error Detail" column (P) on the zh-cn / de-de sheets ---
$ws2.Columns.Item(16).ColumnWidth = 39.17
$ws3.Columns.Item(16).ColumnWidth = 39.17
